$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the mojibake header text in C1 ---
# (was "What is your opinion on the mix of colours and text?Â ", now cleaned up)
$ws.Range("C1").Value = "What is your opinion on the mix of colours and text?"

# --- 2. Remove the stale hidden chart-helper defined names ---
$namesToDelete = @()
foreach ($n in $wb.Names) {
    if ($n.Name -like "_xlchart.v1.*") {
        $namesToDelete += $n.Name
    }
}
foreach ($name in $namesToDelete) {
    $wb.Names.Item($name).Delete()
}

# --- 3. Give column B (the numeric answers to the first question) a fixed width ---
$ws.Columns("B").ColumnWidth = 11.1

# --- 4. New summary row 32: z-test against the neutral midpoint (3) for every
#        numeric question column ---
$ws.Range("A32").Value = "z-test (is neutral (3) or not)"
$ws.Range("B32").Formula = "=_xlfn.Z.TEST(B2:B12,3)"
$ws.Range("C32").Formula = "=_xlfn.Z.TEST(C2:C12,3)"
$ws.Range("D32").Formula = "=_xlfn.Z.TEST(D2:D12,3)"
$ws.Range("G32").Formula = "=_xlfn.Z.TEST(G2:G12,3)"
$ws.Range("H32").Formula = "=_xlfn.Z.TEST(H2:H12,3)"
$ws.Range("I32").Formula = "=_xlfn.Z.TEST(I2:I12,3)"
$ws.Range("K32").Formula = "=_xlfn.Z.TEST(K2:K12,3)"
$ws.Range("N32").Formula = "=_xlfn.Z.TEST(N2:N12,3)"

# --- 5. New summary row 33: t-test comparing the colour+text vs text-only
#        question pairs ---
$ws.Range("A33").Value = "t-test p = (compare text and colour)"
$ws.Range("E33").Formula = "=_xlfn.T.TEST(E2:E12,F2:F12,2,1)"
$ws.Range("L33").Formula = "=_xlfn.T.TEST(L2:L12,M2:M12,2,1)"

# --- 6. Recalculate and match the final selection left behind in the file ---
$excel.Calculate() | Out-Null
$ws.Range("A33").Select() | Out-Null
